$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns remain text, matching the
# original inline-string cell types, even for values that look numeric
# (e.g. "1.003") or percentages.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.821.05"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.639.45"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "309.46"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.3875"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.3808"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").Value = "50.43"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").Value = "1.324"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "0.08383"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "23.66"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "6.962"
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").Value = "7.837"
$ws.Range("E15").Value = "  -3.76%  "
$ws.Range("D16").Value = "0.00001306"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "1.641.71"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "93.56"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "0.06952"
$ws.Range("D20").Value = "19.40"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "6.857"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "13.55"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "23.818.02"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "2.438"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").Value = "2.879"
$ws.Range("E26").Value = "  -8.48%  "
$ws.Range("D27").Value = "21.78"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").Value = "153.25"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "5.548"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("D30").Value = "136.30"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").Value = "7.680"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "2.494"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "1.822.44"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").Value = "0.07984"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").Value = "0.9771"
$ws.Range("E35").Value = "  -6.50%  "
$ws.Range("D36").Value = "0.02891"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("D37").Value = "6.570"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "0.2649"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("D39").Value = "10.41"
$ws.Range("E39").Value = "  -7.75%  "
$ws.Range("D40").Value = "0.09069"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "0.7474"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").Value = "13.25"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").Value = "1.410"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "16.45"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "0.6874"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").Value = "2.411"
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("D47").Value = "4.079"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "0.08208"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "133.89"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Value = "1.212"
$ws.Range("E51").Value = "  -2.27%  "
